$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new daily price records were inserted right after row 144 (existing
# rows 145-263 shift down by two, to 147-265). Insert two blank rows first
# so every following row keeps its original data, then populate the two
# new rows with the new record values.
$ws.Rows("145:146").Insert()

# New row 145: Cebollín "Primera" record for 2021-09-08 (serial 44447)
$ws.Range("A145").Value = 3
$ws.Range("B145").Value = "Femacal de La Calera"
$ws.Range("C145").Value = "Coquimbo"
$ws.Range("D145").Value = 44447
$ws.Range("E145").Value = 5
$ws.Range("F145").Value = 100112037
$ws.Range("G145").Value = "Cebollín"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 270
$ws.Range("K145").Value = 3500
$ws.Range("L145").Value = 4000
$ws.Range("M145").Value = 3722
$ws.Range("N145").Value = "$/paquete 36 unidades"
$ws.Range("O145").Value = "Provincia de Quillota"
$ws.Range("P145").Value = 103
$ws.Range("Q145").Value = 36
$ws.Range("R145").Value = "Hortaliza"

# New row 146: Cebollín "Segunda" record for 2021-09-08 (serial 44447)
$ws.Range("A146").Value = 3
$ws.Range("B146").Value = "Femacal de La Calera"
$ws.Range("C146").Value = "Coquimbo"
$ws.Range("D146").Value = 44447
$ws.Range("E146").Value = 5
$ws.Range("F146").Value = 100112037
$ws.Range("G146").Value = "Cebollín"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Segunda"
$ws.Range("J146").Value = 110
$ws.Range("K146").Value = 2500
$ws.Range("L146").Value = 2500
$ws.Range("M146").Value = 2500
$ws.Range("N146").Value = "$/paquete 36 unidades"
$ws.Range("O146").Value = "Provincia de Quillota"
$ws.Range("P146").Value = 69
$ws.Range("Q146").Value = 36
$ws.Range("R146").Value = "Hortaliza"
